$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextPercent {
    param($cellRef, $text)
    # Assigning a plain "NN%" string via .Value gets auto-parsed by Excel as a
    # percentage number (and forces a new number-format style). To keep the
    # literal text (and the original "General" style), enter it as a formula
    # that evaluates to the text, then convert the formula result to a static
    # value via Copy + PasteSpecial(xlPasteValues).
    $ws.Range($cellRef).Formula = '="' + $text + '"'
    $ws.Range($cellRef).Copy()
    $ws.Range($cellRef).PasteSpecial(-4163)
}

$ws.Range("E2").Value = "2026-02-19 17:48:27"
$ws.Range("E3").Value = "2026-02-19 17:48:29"
$ws.Range("K3").Value = "8.0 MJ/m2"
$ws.Range("E4").Value = "2026-02-19 17:48:32"
$ws.Range("J4").Value = "1009.3 hPa"
$ws.Range("O4").Value = "11.6 °C"
$ws.Range("E5").Value = "2026-02-19 17:48:34"
$ws.Range("I5").Value = "6.6 mm"
$ws.Range("E6").Value = "2026-02-19 17:48:37"
$ws.Range("O6").Value = "10.4 °C"
$ws.Range("E7").Value = "2026-02-19 17:48:39"
$ws.Range("J7").Value = "1010.3 hPa"
$ws.Range("O7").Value = "14.0 °C"
$ws.Range("E8").Value = "2026-02-19 17:48:42"
$ws.Range("J8").Value = "1009.9 hPa"
$ws.Range("K8").Value = "13.4 MJ/m2"
$ws.Range("E9").Value = "2026-02-19 17:48:44"
$ws.Range("K9").Value = "10.6 MJ/m2"
$ws.Range("O9").Value = "10.7 °C"
$ws.Range("E10").Value = "2026-02-19 17:48:47"
$ws.Range("E11").Value = "2026-02-19 17:48:49"
Set-TextPercent "H11" "69%"
$ws.Range("O11").Value = "5.2 °C"
$ws.Range("E12").Value = "2026-02-19 17:48:52"
$ws.Range("E13").Value = "2026-02-19 17:48:54"
Set-TextPercent "H13" "65%"
$ws.Range("J13").Value = "1010.5 hPa"
$ws.Range("K13").Value = "13.7 MJ/m2"
$ws.Range("O13").Value = "4.2 °C"
$ws.Range("E14").Value = "2026-02-19 17:48:56"
$ws.Range("O14").Value = "13.4 °C"
$ws.Range("E15").Value = "2026-02-19 17:48:59"
$ws.Range("E16").Value = "2026-02-19 17:49:01"
$ws.Range("O16").Value = "-7.1 °C"
$ws.Range("E17").Value = "2026-02-19 17:49:04"
$ws.Range("K17").Value = "11.9 MJ/m2"
$ws.Range("E18").Value = "2026-02-19 17:49:06"
Set-TextPercent "H18" "60%"
$ws.Range("J18").Value = "1009.6 hPa"
$ws.Range("O18").Value = "11.7 °C"
$ws.Range("E19").Value = "2026-02-19 17:49:09"
$ws.Range("O19").Value = "5.5 °C"
$ws.Range("E20").Value = "2026-02-19 17:49:11"
$ws.Range("O20").Value = "-5.7 °C"
$ws.Range("E21").Value = "2026-02-19 17:49:14"
Set-TextPercent "H21" "64%"
$ws.Range("J21").Value = "1010.3 hPa"
$ws.Range("O21").Value = "6.3 °C"
$ws.Range("E22").Value = "2026-02-19 17:49:16"
Set-TextPercent "H22" "82%"
$ws.Range("K22").Value = "15.6 MJ/m2"
$ws.Range("E23").Value = "2026-02-19 17:49:18"
$ws.Range("I23").Value = "6.9 mm"
$ws.Range("E24").Value = "2026-02-19 17:49:21"
$ws.Range("J24").Value = "1014.0 hPa"
$ws.Range("E25").Value = "2026-02-19 17:49:23"
$ws.Range("I25").Value = "3.4 mm"
$ws.Range("O25").Value = "-4.3 °C"
$ws.Range("E26").Value = "2026-02-19 17:49:26"
$ws.Range("J26").Value = "1009.3 hPa"
$ws.Range("K26").Value = "9.0 MJ/m2"
$ws.Range("L26").Value = "59.8 km/h - 318º 17:08 TU"
$ws.Range("E27").Value = "2026-02-19 17:49:28"
Set-TextPercent "H27" "68%"
$ws.Range("E28").Value = "2026-02-19 17:49:31"
$ws.Range("O28").Value = "9.1 °C"
$ws.Range("E29").Value = "2026-02-19 17:49:33"
Set-TextPercent "H29" "73%"
$ws.Range("E30").Value = "2026-02-19 17:49:36"
Set-TextPercent "H30" "77%"
$ws.Range("J30").Value = "1009.4 hPa"
$ws.Range("E31").Value = "2026-02-19 17:49:38"
Set-TextPercent "H31" "52%"
$ws.Range("J31").Value = "1008.8 hPa"
$ws.Range("O31").Value = "11.6 °C"
$ws.Range("E32").Value = "2026-02-19 17:49:41"
$ws.Range("E33").Value = "2026-02-19 17:49:43"
Set-TextPercent "H33" "61%"
$ws.Range("O33").Value = "3.5 °C"
$ws.Range("E34").Value = "2026-02-19 17:49:46"
$ws.Range("E35").Value = "2026-02-19 17:49:48"
Set-TextPercent "H35" "68%"
$ws.Range("J35").Value = "1015.5 hPa"
$ws.Range("E36").Value = "2026-02-19 17:49:51"
$ws.Range("J36").Value = "1009.7 hPa"
$ws.Range("O36").Value = "12.0 °C"
$ws.Range("E37").Value = "2026-02-19 17:49:53"
Set-TextPercent "H37" "73%"
$ws.Range("O37").Value = "5.6 °C"
$ws.Range("E38").Value = "2026-02-19 17:49:55"
Set-TextPercent "H38" "58%"
$ws.Range("E39").Value = "2026-02-19 17:49:58"
$ws.Range("L39").Value = "78.8 km/h - 293º 17:15 TU"
$ws.Range("E40").Value = "2026-02-19 17:50:00"
Set-TextPercent "H40" "76%"
$ws.Range("O40").Value = "6.1 °C"
$ws.Range("E41").Value = "2026-02-19 17:50:02"
$ws.Range("J41").Value = "1012.2 hPa"
$ws.Range("E42").Value = "2026-02-19 17:50:05"
$ws.Range("E43").Value = "2026-02-19 17:50:07"
$ws.Range("O43").Value = "9.2 °C"
$ws.Range("E44").Value = "2026-02-19 17:50:10"
$ws.Range("I44").Value = "7.6 mm"
$ws.Range("K44").Value = "10.2 MJ/m2"
$ws.Range("E45").Value = "2026-02-19 17:50:12"
Set-TextPercent "H45" "86%"
$ws.Range("J45").Value = "1014.5 hPa"
$ws.Range("E46").Value = "2026-02-19 17:50:14"
$ws.Range("J46").Value = "1014.9 hPa"

$excel.CutCopyMode = 0

